# Vaccine surveillance report: add week 5 death28 data + render format
# (see commit message: "add week 5 data and render format")

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "diag_week_5" (sheet4): the header row's wrapped text needs more
#    vertical room now that the death sheets below grew taller too.
# ---------------------------------------------------------------------
$wsDiag5 = $wb.Worksheets.Item("diag_week_5")
$wsDiag5.Rows.Item(1).RowHeight = 97.2

# ---------------------------------------------------------------------
# 2. "death28_week_5" (sheet5): this sheet still had a stray leftover
#    footnote row ("[1]" marker in G2/H2) instead of the real "Under 18"
#    data row. Delete that placeholder row so the real age-band rows
#    (previously rows 3-10) shift up into rows 2-9, matching the other
#    completed week-5 sheets.
# ---------------------------------------------------------------------
$wsDeath28_5 = $wb.Worksheets.Item("death28_week_5")
$wsDeath28_5.Rows.Item(2).Delete()

# Reflect the new selection left behind on this sheet (row 2 is now the
# first real data row, so the whole row is the live selection).
$wsDeath28_5.Activate()
$wsDeath28_5.Range("A2:XFD2").Select()

# ---------------------------------------------------------------------
# 3. "death_week_5" (sheet6): becomes the active/selected sheet/tab now
#    that the week-5 data entry is finished.
# ---------------------------------------------------------------------
$wsDeath_5 = $wb.Worksheets.Item("death_week_5")
$wsDeath_5.Activate()
$wsDeath_5.Range("B1").Select()

# Best-effort: restore the saved window size recorded the last time the
# workbook was closed (cosmetic only; not exposed by every host).
try {
    $win = $excel.ActiveWindow
    $win.Width = 23040
    $win.Height = 6732
} catch {}
